$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 currently holds the combined Cleric/Rogue feature strings.
# Split them: keep the "final tier" portion in row 9, and move the
# earlier tiers into a brand-new row 10 (copying row 9's formatting).
$ws.Range("A9:B9").Copy()
$ws.Range("A10:B10").PasteSpecial(-4122)

$ws.Range("A9").Value = "10/Divine Intervention"
$ws.Range("B9").Value = "1/Thieves' Cant=1/Expertise=2/Cunning Action=5/Uncanny Dodge=7/Evasion=11/Reliable Talent=14/Blindsense=15/Slippery Mind=18/Elusive=20/Stroke of Luck"

$ws.Range("A10").Value = "2/Channel Divinity=5/Destroy Undead"
$ws.Range("B10").Value = "1/Sneak Attack"
